# Apply the sweep-measurements update: refresh row 2 of the "Measurements"
# sheet with the new run's captured values, and mirror those same values
# into the corresponding Max/Min/Mean rows of the "Statistics" sheet
# (there is a single test in this workbook, so Max == Min == Mean == the
# row-2 value for every numeric metric).

$wb = $excel.ActiveWorkbook
$measurements = $wb.Worksheets.Item("Measurements")
$statistics = $wb.Worksheets.Item("Statistics")

# ---------------------------------------------------------------------
# 1) Measurements!row2 - new captured values
# ---------------------------------------------------------------------
$measurements.Range("A2").Value = 1.230139493942261
$measurements.Range("B2").Value = 16.80768704414368
# C2, D2, E2 unchanged (2, 6, 1)
$measurements.Range("F2").Value = 0.6845536231994629
$measurements.Range("G2").Value = -12.03915
$measurements.Range("H2").Value = 6.392822
$measurements.Range("I2").Value = 6.047246933
$measurements.Range("J2").Value = -45.75321198
$measurements.Range("K2").Value = 0.838693380355835
$measurements.Range("L2").Value = 6.15513134003
$measurements.Range("M2").Value = -51.4895620346
$measurements.Range("N2").Value = -52.3492879868
$measurements.Range("O2").Value = 0.6965670585632324
$measurements.Range("P2").Value = 6.083156586
$measurements.Range("Q2").Value = -51.17642212
$measurements.Range("R2").Value = 0.9600005149841309
$measurements.Range("S2").Value = 6.19470405579
$measurements.Range("T2").Value = -52.3015956879
$measurements.Range("U2").Value = -53.7921848297
$measurements.Range("V2").Value = 0.6848137378692627
$measurements.Range("W2").Value = 6.300565719604492
$measurements.Range("X2").Value = 2
$measurements.Range("Y2").Value = 1.128886699676514
# Z2 unchanged (0)

$measurements.Range("AW2").Value = 9.385

$measurements.Range("BC2").Value = "-45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75, -45.75"
$measurements.Range("BD2").Value = 0.862
$measurements.Range("BF2").Value = "-51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18, -51.18"
$measurements.Range("BG2").Value = 0.862

# ---------------------------------------------------------------------
# 2) Statistics sheet - Max / Min / Mean rows mirror the single test run
# ---------------------------------------------------------------------
$statistics.Range("B3").Value = 1.230139493942261
$statistics.Range("B4").Value = 1.230139493942261
$statistics.Range("B5").Value = 1.230139493942261

$statistics.Range("B6").Value = 16.80768704414368
$statistics.Range("B7").Value = 16.80768704414368
$statistics.Range("B8").Value = 16.80768704414368

$statistics.Range("B18").Value = 0.6845536231994629
$statistics.Range("B19").Value = 0.6845536231994629
$statistics.Range("B20").Value = 0.6845536231994629

$statistics.Range("B21").Value = -12.03915
$statistics.Range("B22").Value = -12.03915
$statistics.Range("B23").Value = -12.03915

$statistics.Range("B24").Value = 6.392822
$statistics.Range("B25").Value = 6.392822
$statistics.Range("B26").Value = 6.392822

$statistics.Range("B27").Value = 6.047246933
$statistics.Range("B28").Value = 6.047246933
$statistics.Range("B29").Value = 6.047246933

$statistics.Range("B30").Value = -45.75321198
$statistics.Range("B31").Value = -45.75321198
$statistics.Range("B32").Value = -45.75321198

$statistics.Range("B33").Value = 0.838693380355835
$statistics.Range("B34").Value = 0.838693380355835
$statistics.Range("B35").Value = 0.838693380355835

$statistics.Range("B36").Value = 6.15513134003
$statistics.Range("B37").Value = 6.15513134003
$statistics.Range("B38").Value = 6.15513134003

$statistics.Range("B39").Value = -51.4895620346
$statistics.Range("B40").Value = -51.4895620346
$statistics.Range("B41").Value = -51.4895620346

$statistics.Range("B42").Value = -52.3492879868
$statistics.Range("B43").Value = -52.3492879868
$statistics.Range("B44").Value = -52.3492879868

$statistics.Range("B45").Value = 0.6965670585632324
$statistics.Range("B46").Value = 0.6965670585632324
$statistics.Range("B47").Value = 0.6965670585632324

$statistics.Range("B48").Value = 6.083156586
$statistics.Range("B49").Value = 6.083156586
$statistics.Range("B50").Value = 6.083156586

$statistics.Range("B51").Value = -51.17642212
$statistics.Range("B52").Value = -51.17642212
$statistics.Range("B53").Value = -51.17642212

$statistics.Range("B54").Value = 0.9600005149841309
$statistics.Range("B55").Value = 0.9600005149841309
$statistics.Range("B56").Value = 0.9600005149841309

$statistics.Range("B57").Value = 6.19470405579
$statistics.Range("B58").Value = 6.19470405579
$statistics.Range("B59").Value = 6.19470405579

$statistics.Range("B60").Value = -52.3015956879
$statistics.Range("B61").Value = -52.3015956879
$statistics.Range("B62").Value = -52.3015956879

$statistics.Range("B63").Value = -53.7921848297
$statistics.Range("B64").Value = -53.7921848297
$statistics.Range("B65").Value = -53.7921848297

$statistics.Range("B66").Value = 0.6848137378692627
$statistics.Range("B67").Value = 0.6848137378692627
$statistics.Range("B68").Value = 0.6848137378692627

$statistics.Range("B69").Value = 6.300565719604492
$statistics.Range("B70").Value = 6.300565719604492
$statistics.Range("B71").Value = 6.300565719604492

$statistics.Range("B72").Value = 2
$statistics.Range("B73").Value = 2
$statistics.Range("B74").Value = 2

$statistics.Range("B75").Value = 1.128886699676514
$statistics.Range("B76").Value = 1.128886699676514
$statistics.Range("B77").Value = 1.128886699676514

$statistics.Range("B81").Value = 9.385
$statistics.Range("B82").Value = 9.385
$statistics.Range("B83").Value = 9.385

$statistics.Range("B93").Value = 0.862
$statistics.Range("B94").Value = 0.862
$statistics.Range("B95").Value = 0.862

$statistics.Range("B96").Value = 0.862
$statistics.Range("B97").Value = 0.862
$statistics.Range("B98").Value = 0.862
